# Adjust property of scene: update the camera offset position/rotation
# values (CamOffestPos / CamOffestRot, columns J and K) for the "Demo1"
# row (row 2) and the "SelectScene" row (row 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Demo1 (row 2)
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# SelectScene (row 6)
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Leave the selection on the last-edited cell, matching the author's
# recorded cursor position when the workbook was saved.
$ws.Range("K7").Select()
